$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 176.38461
$ws.Range("I2").Value = 98
$ws.Range("J2").Value = 199.9
$ws.Range("K2").Value = 98
$ws.Range("L2").Value = 199.9
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = -425.9
$ws.Range("H3").Value = 42139.8
$ws.Range("J3").Value = 42139.8
$ws.Range("L3").Value = 42139.8
$ws.Range("N3").Value = -42367.8
$ws.Range("H42").Value = 1306.3572
$ws.Range("I42").Value = 579.2
$ws.Range("J42").Value = 3124.25
$ws.Range("K42").Value = 1737.6
$ws.Range("L42").Value = 9372.75
$ws.Range("M42").Value = -1507.6
$ws.Range("N42").Value = -9832.75
$ws.Range("H86").Value = 5249.5713
$ws.Range("J86").Value = 5999.4
$ws.Range("L86").Value = 5999.4
$ws.Range("N86").Value = -8245.4
$ws.Range("H89").Value = 5249.5713
$ws.Range("J89").Value = 5999.4
$ws.Range("L89").Value = 29997
$ws.Range("N89").Value = -41229
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H102").Value = 42139.8
$ws.Range("J102").Value = 42139.8
$ws.Range("L102").Value = 42139.8
$ws.Range("N102").Value = -48629.8
$ws.Range("H113").Value = 7461.2285
$ws.Range("J113").Value = 8185.533
$ws.Range("L113").Value = 8185.533
$ws.Range("N113").Value = -14693.533
$ws.Range("H125").Value = 2859.6667
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 18000
$ws.Range("M125").Value = -15540
$ws.Range("H132").Value = 4383.619
$ws.Range("I132").Value = 4195.1797
$ws.Range("J132").Value = 6833.3335
$ws.Range("K132").Value = 12585.5391
$ws.Range("L132").Value = 20500.0005
$ws.Range("M132").Value = -10055.5391
$ws.Range("N132").Value = -25560.0005
$ws.Range("H135").Value = 1370.5
$ws.Range("I135").Value = 1224.1666
$ws.Range("K135").Value = 11017.4994
$ws.Range("M135").Value = -8482.499400000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4833.0435
$ws.Range("I74").Value = 2750.25
$ws.Range("J74").Value = 5943.8667
$ws.Range("K74").Value = 2750.25
$ws.Range("L74").Value = 5943.8667
$ws.Range("M74").Value = -1876.25
$ws.Range("N74").Value = -7691.8667
$ws.Range("H77").Value = 4833.0435
$ws.Range("I77").Value = 2750.25
$ws.Range("J77").Value = 5943.8667
$ws.Range("K77").Value = 13751.25
$ws.Range("L77").Value = 29719.3335
$ws.Range("M77").Value = -9383.25
$ws.Range("N77").Value = -38455.33349999999
$ws.Range("H88").Value = 13103.777
$ws.Range("J88").Value = 18661.334
$ws.Range("L88").Value = 18661.334
$ws.Range("N88").Value = -19473.334
$ws.Range("H91").Value = 13103.777
$ws.Range("J91").Value = 18661.334
$ws.Range("L91").Value = 18661.334
$ws.Range("N91").Value = -21469.334
$ws.Range("H101").Value = 57496.75
$ws.Range("J101").Value = 57496.75
$ws.Range("L101").Value = 57496.75
$ws.Range("N101").Value = -63986.75
$ws.Range("H122").Value = 3089.9167
$ws.Range("I122").Value = 3143.5454
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 9430.636200000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6980.636200000001
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 1925.7142
$ws.Range("I132").Value = 1498.5
$ws.Range("K132").Value = 4495.5
$ws.Range("M132").Value = -1965.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7326.148
$ws.Range("I86").Value = 2992.389
$ws.Range("K86").Value = 2992.389
$ws.Range("M86").Value = -1869.389
$ws.Range("H89").Value = 7326.148
$ws.Range("I89").Value = 2992.389
$ws.Range("K89").Value = 14961.945
$ws.Range("M89").Value = -9345.945
$ws.Range("H99").Value = 4923
$ws.Range("J99").Value = 5000
$ws.Range("L99").Value = 5000
$ws.Range("N99").Value = -7996

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 617633.6
$ws.Range("I6").Value = 704438.4399999999
$ws.Range("K6").Value = 704438.4399999999
$ws.Range("M6").Value = -704325.4399999999
$ws.Range("H29").Value = 899
$ws.Range("J29").Value = 899
$ws.Range("L29").Value = 899
$ws.Range("N29").Value = -1485
$ws.Range("H132").Value = 3001.2173
$ws.Range("I132").Value = 3237.9285
$ws.Range("J132").Value = 2633
$ws.Range("K132").Value = 9713.7855
$ws.Range("L132").Value = 7899
$ws.Range("M132").Value = -7183.7855
$ws.Range("N132").Value = -12959
$ws.Range("H134").Value = 4957.45
$ws.Range("I134").Value = 4472.515
$ws.Range("K134").Value = 13417.545
$ws.Range("M134").Value = -10882.545

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 17857.143
$ws.Range("J93").Value = 17857.143
$ws.Range("L93").Value = 53571.429
$ws.Range("N93").Value = -57315.429
$ws.Range("H131").Value = 17096098
$ws.Range("J131").Value = 22224790
$ws.Range("L131").Value = 66674370
$ws.Range("N131").Value = -66684450

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3437.182
$ws.Range("I80").Value = 4509.5
$ws.Range("J80").Value = 2150.4
$ws.Range("K80").Value = 4509.5
$ws.Range("L80").Value = 2150.4
$ws.Range("M80").Value = -3511.5
$ws.Range("N80").Value = -4146.4
$ws.Range("H83").Value = 3437.182
$ws.Range("I83").Value = 4509.5
$ws.Range("J83").Value = 2150.4
$ws.Range("K83").Value = 22547.5
$ws.Range("L83").Value = 10752
$ws.Range("M83").Value = -17555.5
$ws.Range("N83").Value = -20736
$ws.Range("H97").Value = 1239.8
$ws.Range("I97").Value = 924.75
$ws.Range("J97").Value = 2500
$ws.Range("K97").Value = 924.75
$ws.Range("L97").Value = 2500
$ws.Range("M97").Value = -428.75
$ws.Range("N97").Value = -3492
$ws.Range("H102").Value = 4961.1
$ws.Range("I102").Value = 5068
$ws.Range("K102").Value = 5068
$ws.Range("M102").Value = -3446
$ws.Range("H126").Value = 6134.6665
$ws.Range("I126").Value = 2570.6667
$ws.Range("K126").Value = 7712.000100000001
$ws.Range("M126").Value = -5242.000100000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1975
$ws.Range("I100").Value = 1950
$ws.Range("K100").Value = 1950
$ws.Range("M100").Value = -1409
$ws.Range("H122").Value = 7057.2856
$ws.Range("I122").Value = 4499
$ws.Range("J122").Value = 7755
$ws.Range("K122").Value = 13497
$ws.Range("L122").Value = 23265
$ws.Range("M122").Value = -11047
$ws.Range("N122").Value = -28165
$ws.Range("H132").Value = 3445.5
$ws.Range("I132").Value = 2418.375
$ws.Range("J132").Value = 4472.625
$ws.Range("K132").Value = 7255.125
$ws.Range("L132").Value = 13417.875
$ws.Range("M132").Value = -4725.125
$ws.Range("N132").Value = -18477.875
$ws.Range("H136").Value = 2953.1943
$ws.Range("I136").Value = 1273.3
$ws.Range("J136").Value = 3599.3076
$ws.Range("K136").Value = 3819.9
$ws.Range("L136").Value = 10797.9228
$ws.Range("M136").Value = -1269.9
$ws.Range("N136").Value = -15897.9228

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4750
$ws.Range("I122").Value = 4612
$ws.Range("J122").Value = 4842
$ws.Range("K122").Value = 13836
$ws.Range("L122").Value = 14526
$ws.Range("M122").Value = -11386
$ws.Range("N122").Value = -19426
$ws.Range("H132").Value = 2720.1428
$ws.Range("I132").Value = 2363.6667
$ws.Range("K132").Value = 7091.000100000001
$ws.Range("M132").Value = -4561.000100000001
$ws.Range("H136").Value = 8675.521000000001
$ws.Range("I136").Value = 9576.6
$ws.Range("K136").Value = 28729.8
$ws.Range("M136").Value = -26179.8
